# Applies the "CVD files" update described in the commit message:
# a handful of data corrections across existing location sheets, plus
# two brand-new location sheets (Betzdorf Germany, Cotia São Paulo Brazil)
# appended at the end of the workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Baja California Mexico
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Baja California Mexico")
$ws.Range("L4").Value = 0.0161
$ws.Range("L5").ClearContents()
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("I8").Value = 0.0034
$ws.Range("L8:W8").Value = 0

# ---------------------------------------------------------------------
# 2. Fremont California
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fremont California")
$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776

# ---------------------------------------------------------------------
# 3. Kristianstad Sweden
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Kristianstad Sweden")
$ws.Range("L7").ClearContents()
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776
$ws.Range("L10").Value = 0.0105
$ws.Range("M10:W10").Value = 0

# ---------------------------------------------------------------------
# 4. Marengo Illinois
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Marengo Illinois")
$ws.Range("L4").Value = 0.0182
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("G9").Value = 0.0255
$ws.Range("H9").Value = 0.0064
$ws.Range("J9").Value = 0.0322
$ws.Range("K9").Value = 0.0133
$ws.Range("L9").Value = 0.0135

# ---------------------------------------------------------------------
# 5. Milwaukee Pmc Hq Wisconsin
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("L5").ClearContents()

# ---------------------------------------------------------------------
# 6. Rock Road Radford Virginia
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Rock Road Radford Virginia")
$ws.Range("L4").Value = 0.0161
$ws.Range("L7").ClearContents()

# ---------------------------------------------------------------------
# 7. Sandy Point Town St Kitts
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sandy Point Town St Kitts")
$ws.Range("E5").Value = 0.0776
$ws.Range("E6").Value = 0.0776
$ws.Range("E7").Value = 0.0776

# ---------------------------------------------------------------------
# 8. Tianjin China
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tianjin China")
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

# ---------------------------------------------------------------------
# 9. Wolfschlugen Germany
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Wolfschlugen Germany")
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776

# ---------------------------------------------------------------------
# 10. Bristol Connecticut
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Bristol Connecticut")
$ws.Range("E4").Value = 0.0776
$ws.Range("E5").Value = 0.0776

# ---------------------------------------------------------------------
# 11. Cleveland Ohio
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Cleveland Ohio")
$ws.Range("E7").Value = 0.0776
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("K9").Value = 0.0222
$ws.Range("L9").Value = 0.0108

# ---------------------------------------------------------------------
# 12. Devon United Kingdom
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Devon United Kingdom")
$ws.Range("L7").ClearContents()

# ---------------------------------------------------------------------
# 13. Downers Grove Illinois
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Downers Grove Illinois")
$ws.Range("I4").Value = 0.0345
$ws.Range("J4").Value = 0.0333
$ws.Range("L7").Value = 0.6667

# ---------------------------------------------------------------------
# 14. East Aurora New York
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("East Aurora New York")
$ws.Range("E5").Value = 0.571428571428571
$ws.Range("E6").Value = 0.571428571428571
$ws.Range("E7").Value = 0.571428571428571
$ws.Range("L7").Value = 1
$ws.Range("M7:W7").Value = 0.571428571428571
$ws.Range("E8").Value = 0.0776
$ws.Range("E9").Value = 0.0776
$ws.Range("E10").Value = 0.0776
$ws.Range("G10").Value = 0.0413
$ws.Range("H10").Value = 0.008
$ws.Range("I10").Value = 0.0397
$ws.Range("J10").Value = 0.0887
$ws.Range("K10").Value = 0.0079
$ws.Range("L10").Value = 0.0397

# ---------------------------------------------------------------------
# 15. New sheet: Betzdorf Germany (appended after the last sheet)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$betzdorf = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$betzdorf.Name = "Betzdorf Germany"

$headers = @("segment_function","division_function","location","cvd","ytd","data_source","Jan","Feb","Mar","Q1","Apr","May","Jun","Q2","Jul","Aug","Sep","Q3","Oct","Nov","Dec","Q4","FY")
$betzdorf.Range("A1:W1").Value = $headers

$betzdorf.Range("E2:W3").NumberFormat = "0.0%"

$betzdorf.Range("A2").Value = "AMC"
$betzdorf.Range("B2").Value = "AMC Linear Motion Division"
$betzdorf.Range("C2").Value = "Betzdorf Germany"
$betzdorf.Range("D2").Value = "Professional Voluntary Turnover"
$betzdorf.Range("E2").Value = 0
$betzdorf.Range("F2").Value = "Commit/Forecast"
$betzdorf.Range("G2:W2").Value = 0

$betzdorf.Range("A3").Value = "AMC"
$betzdorf.Range("B3").Value = "AMC Linear Motion Division"
$betzdorf.Range("C3").Value = "Betzdorf Germany"
$betzdorf.Range("D3").Value = "Internal Fill Rate"
$betzdorf.Range("E3").Value = 1
$betzdorf.Range("F3").Value = "Commit/Forecast"
$betzdorf.Range("K3").Value = 1
$betzdorf.Range("M3:W3").Value = 1

# ---------------------------------------------------------------------
# 16. New sheet: Cotia São Paulo Brazil (appended after Betzdorf Germany)
# ---------------------------------------------------------------------
$cotia = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $betzdorf)
$cotia.Name = "Cotia São Paulo Brazil"

$cotia.Range("A1:W1").Value = $headers

$cotia.Range("E2:W3").NumberFormat = "0.0%"

$cotia.Range("A2").Value = "AMC"
$cotia.Range("B2").Value = "AMC Linear Motion Division"
$cotia.Range("C2").Value = "Cotia São Paulo Brazil"
$cotia.Range("D2").Value = "Professional Voluntary Turnover"
$cotia.Range("E2").Value = 0
$cotia.Range("F2").Value = "Commit/Forecast"
$cotia.Range("G2:W2").Value = 0

$cotia.Range("A3").Value = "AMC"
$cotia.Range("B3").Value = "AMC Linear Motion Division"
$cotia.Range("C3").Value = "Cotia São Paulo Brazil"
$cotia.Range("D3").Value = "Manufacturing Voluntary Turnover"
$cotia.Range("E3").Value = 0.0776
$cotia.Range("F3").Value = "Commit/Forecast"
$cotia.Range("G3:W3").Value = 0
